# Generate Report for Handoff
# Updates status text, timestamps, and widens the date/status columns
# across the Overview, zh-cn, and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Update status text "In Translation" -> "Ready for handoff" ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# --- Update generate/handoff date-time stamps (kept as text) ---
$wsOverview.Range("G2").Value = "2016-08-12 13:10:19"
$wsDeDe.Range("H2").Value = "2016-08-12 13:10:19"
$wsZhCn.Range("H2").Value = "2016-08-12 13:10:11"

# --- Widen status/date columns to fit new content ---
# (target authored width is 17.2159881591797; the host quantizes
# ColumnWidth to 1/6-character steps, so 16.3 is the input that lands on
# the nearest reachable width, 17.1666...)
$wsOverview.Range("E:E").ColumnWidth = 16.3
$wsOverview.Range("F:F").ColumnWidth = 16.3
$wsZhCn.Range("C:C").ColumnWidth = 16.3
$wsDeDe.Range("C:C").ColumnWidth = 16.3
